# "added 4wk low sales check"
# Updates the forecast figures (MyForecast, Inventory Coverage, Stockout Risk,
# Reorder Urgency, Seasonality Index) on the "Forecast Comparison" sheet and
# the dependent rollup figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$forecast = $wb.Worksheets.Item("Forecast Comparison")
$summary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet --------------------------------------------
# Row 2 (W10)
$forecast.Range("D2").Value = 107
$forecast.Range("H2").Value = 5.91
$forecast.Range("L2").Value = 1.12

# Row 3 (W11)
$forecast.Range("D3").Value = 132
$forecast.Range("H3").Value = 3.98
$forecast.Range("L3").Value = 1.03

# Row 4 (W12)
$forecast.Range("D4").Value = 152
$forecast.Range("H4").Value = 2.59
$forecast.Range("L4").Value = 0.96

# Row 5 (W13)
$forecast.Range("D5").Value = 154
$forecast.Range("H5").Value = 1.56
$forecast.Range("L5").Value = 0.97

# Row 6 (W14)
$forecast.Range("D6").Value = 135
$forecast.Range("H6").Value = 0.64
$forecast.Range("J6").Value = "Urgent"
$forecast.Range("L6").Value = 0.99

# Row 7 (W15)
$forecast.Range("D7").Value = 108
$forecast.Range("H7").Value = 0
$forecast.Range("I7").Value = "High"
$forecast.Range("J7").Value = "Urgent"
$forecast.Range("L7").Value = 1.01

# Row 8 (W16)
$forecast.Range("D8").Value = 89
$forecast.Range("H8").Value = 0
$forecast.Range("I8").Value = "High"
$forecast.Range("J8").Value = "Urgent"
$forecast.Range("L8").Value = 0.9

# Row 9 (W17)
$forecast.Range("D9").Value = 91
$forecast.Range("H9").Value = 0
$forecast.Range("L9").Value = 0.95

# Row 10 (W18)
$forecast.Range("D10").Value = 106
$forecast.Range("L10").Value = 1.04

# Row 11 (W19)
$forecast.Range("D11").Value = 114
$forecast.Range("L11").Value = 1.07

# Row 12 (W20)
$forecast.Range("D12").Value = 101
$forecast.Range("L12").Value = 1.07

# Row 13 (W21)
$forecast.Range("D13").Value = 69
$forecast.Range("L13").Value = 0.84

# Row 14 (W22)
$forecast.Range("D14").Value = 36
$forecast.Range("L14").Value = 1.09

# Row 15 (W23)
$forecast.Range("D15").Value = 27
$forecast.Range("L15").Value = 0.93

# Row 16 (W24)
$forecast.Range("D16").Value = 37
$forecast.Range("L16").Value = 1.18

# Row 17 (W25)
$forecast.Range("D17").Value = 57
$forecast.Range("L17").Value = 1.17

# --- Summary sheet ----------------------------------------------------------
# These "numbers" are stored as text on this sheet (same as all other cells
# here), so a leading apostrophe is used to keep Excel from re-typing them
# as numeric values.
$summary.Range("B9").Value  = "'1515"
$summary.Range("B10").Value = "'968"
$summary.Range("B11").Value = "'545"
$summary.Range("B12").Value = "'154"
$summary.Range("B14").Value = "'27"
